# Commit: "create new slide 3"
#
# The deck currently has 2 slides:
#   Slide 1 - Title slide (layout "Титульный слайд")
#   Slide 2 - Title + Content slide (layout "Заголовок и объект"), with
#             empty title ("Заголовок 1") and content ("Содержимое 2")
#             placeholders.
#
# The target adds a brand-new slide 3 that is structurally identical to
# slide 2 (same layout, same two empty placeholders). The simplest and
# most faithful way to reproduce that is to duplicate slide 2 and let the
# duplicate land right after it, as slide 3.

$p = $ppt.ActivePresentation

$source = $p.Slides.Item(2)
$source.Duplicate() | Out-Null
